$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 2) re-population: the "NAME" / "ACCOUNT NO" / "DATE
# RELEASED" / "LOAN TERM" / "AMOUNT RELEASED" columns (C:G) were reshuffled
# into: ACCOUNT NO, NAME, DATE RELEASED, AMOUNT RELEASED, LOAN TERM.
# Both the header text AND its cell formatting (date / accounting number
# format) travel together, so we round-trip through staging cells far off
# to the right (AA:AE) using copy/paste-special so the existing styles
# (General / date / accounting) are preserved exactly rather than minted
# as brand new styles.
# ---------------------------------------------------------------------------

$src = @("C2", "D2", "E2", "F2", "G2")
$stage = @("AA2", "AB2", "AC2", "AD2", "AE2")

for ($i = 0; $i -lt $src.Length; $i++) {
    $ws.Range($src[$i]).Copy()
    $ws.Range($stage[$i]).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($src[$i]).Copy()
    $ws.Range($stage[$i]).PasteSpecial(-4163)  # xlPasteValues
}

# new C2..G2 <- old E2, F2, C2, G2, D2
$order = @(2, 3, 0, 4, 1)
for ($i = 0; $i -lt $src.Length; $i++) {
    $from = $stage[$order[$i]]
    $to = $src[$i]
    $ws.Range($from).Copy()
    $ws.Range($to).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($from).Copy()
    $ws.Range($to).PasteSpecial(-4163)  # xlPasteValues
}

# clean up the staging cells
$ws.Range("AA2:AE2").Clear()

# ---------------------------------------------------------------------------
# Selection / active cell moved to F18 while populating the form.
# ---------------------------------------------------------------------------
$null = $ws.Range("F18").Select()

# ---------------------------------------------------------------------------
# Explicit portrait page setup (as saved by the editing client).
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
